$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume list refresh (GitHub Actions scheduled update)

# Row 2
$ws.Range("D2").Value = '23.479.39'
$ws.Range("E2").Value = '  +1.25%  '

# Row 3
$ws.Range("D3").Value = '1.637.54'
$ws.Range("E3").Value = '  +2.22%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.11'
$ws.Range("E5").Value = '  +1.38%  '

# Row 6
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.08%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3768'
$ws.Range("E7").Value = '  -0.38%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.30'
$ws.Range("E8").Value = '  +1.15%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3647'
$ws.Range("E9").Value = '  +0.97%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.270'
$ws.Range("E10").Value = '  +0.53%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08170'
$ws.Range("E11").Value = '  +0.53%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.01%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.03'
$ws.Range("E13").Value = '  +1.80%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.643'
$ws.Range("E14").Value = '  +0.98%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001278'
$ws.Range("E15").Value = '  +2.43%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.380'
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").Value = '1.636.38'
$ws.Range("E17").Value = '  +2.16%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.82'
$ws.Range("E18").Value = '  +1.05%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06970'
$ws.Range("E19").Value = '  +1.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.22'
$ws.Range("E20").Value = '  +0.93%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.546'
$ws.Range("E21").Value = '  +0.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.00%  '

# Row 23
$ws.Range("D23").Value = '23.512.82'
$ws.Range("E23").Value = '  +1.43%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.81'
$ws.Range("E24").Value = '  -1.00%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.101'
$ws.Range("E25").Value = '  +4.25%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.417'
$ws.Range("E26").Value = '  +1.30%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.28'

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.35'
$ws.Range("E28").Value = '  +1.34%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.364'
$ws.Range("E29").Value = '  +2.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.52'
$ws.Range("E30").Value = '  +1.32%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.355'
$ws.Range("E31").Value = '  -0.98%  '

# Row 32
$ws.Range("D32").Value = '1.819.06'
$ws.Range("E32").Value = '  +2.25%  '

# Row 33
$ws.Range("E33").Value = '  -0.58%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9656'
$ws.Range("E34").Value = '  -0.59%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02828'
$ws.Range("E35").Value = '  +4.05%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.34'
$ws.Range("E36").Value = '  +0.23%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07366'
$ws.Range("E37").Value = '  -2.08%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2538'
$ws.Range("E38").Value = '  +1.18%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.177'
$ws.Range("E39").Value = '  +1.04%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08857'
$ws.Range("E40").Value = '  +0.63%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.382'
$ws.Range("E41").Value = '  +1.51%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7124'
$ws.Range("E42").Value = '  +0.58%  '

# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.28'
$ws.Range("E43").Value = '  +4.73%  '

# Row 44
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.52'
$ws.Range("E44").Value = '  +0.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6558'
$ws.Range("E45").Value = '  +0.45%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.341'
$ws.Range("E46").Value = '  +1.39%  '

# Row 47
$ws.Range("E47").Value = '  +0.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.035'
$ws.Range("E48").Value = '  +0.62%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07970'
$ws.Range("E49").Value = '  +0.13%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '129.54'
$ws.Range("E50").Value = '  -1.89%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.208'
$ws.Range("E51").Value = '  +0.43%  '
